$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Computing Time (sec)" -> "Computing Time (ns)" (columns C and F)
$ws.Range("C1").Value = "Computing Time (ns)"
$ws.Range("F1").Value = "Computing Time (ns)"

# Row 2 (f1)
$ws.Range("B2").Value = 159
$ws.Range("C2").Value = 72900
$ws.Range("D2").Value = 53.9
$ws.Range("F2").Value = 82000

# Row 3 (f2)
$ws.Range("B3").Value = 851
$ws.Range("C3").Value = 92700
$ws.Range("D3").Value = 83.11
$ws.Range("F3").Value = 127100

# Row 4 (f3)
$ws.Range("B4").Value = 28
$ws.Range("C4").Value = 53500
$ws.Range("D4").Value = 80
$ws.Range("F4").Value = 56800

# Row 5 (f4)
$ws.Range("C5").Value = 51500
$ws.Range("F5").Value = 52500

# Row 6 (f6)
$ws.Range("B6").Value = 33
$ws.Range("C6").Value = 60100
$ws.Range("D6").Value = 63.46
$ws.Range("F6").Value = 72000

# Row 7 (f7)
$ws.Range("B7").Value = 70
$ws.Range("C7").Value = 53400
$ws.Range("D7").Value = 65.42
$ws.Range("F7").Value = 55200

# Row 8 (f8)
$ws.Range("B8").Value = 9134
$ws.Range("C8").Value = 89100
$ws.Range("D8").Value = 93.52
$ws.Range("F8").Value = 113900

# Row 9 (f9)
$ws.Range("B9").Value = 72
$ws.Range("C9").Value = 55600
$ws.Range("D9").Value = 55.38
$ws.Range("F9").Value = 56200

# Row 10 (f10)
$ws.Range("B10").Value = 850
$ws.Range("C10").Value = 90200
$ws.Range("D10").Value = 82.93000000000001
$ws.Range("F10").Value = 123300
